$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph formatting updates for the four bullet paragraphs under
#    "General provisions" (numId=4): ind 426 -> 425 (both left & hanging),
#    spacing lineRule atLeast -> auto, and removing explicit "before" spacing
#    on the first one.
# ---------------------------------------------------------------------------
foreach ($idx in 6,7,8,9) {
    $p = $d.Paragraphs($idx)
    $p.LineSpacingRule = 5   # wdLineSpaceMultiple -> produces lineRule="auto"
    $p.LineSpacing = 18      # 18pt * 20 = 360 twips, matches w:line="360"
    $p.LeftIndent = 21.25    # 425 twips
    $p.FirstLineIndent = -21.25
    $p.SpaceAfter = 0
}
# Paragraph 6 ("The data administrator is ...") also drops its "space before".
$d.Paragraphs(6).SpaceBefore = 0

# ---------------------------------------------------------------------------
# 2) Merge the split runs around "with regard to" into a single run (removes
#    the now redundant proofErr gramStart/gramEnd wrapper).
# ---------------------------------------------------------------------------
$oldText1 = "individuals with regard to the processing of personal data and on the free movement of such data, and repealing Directive 95/46 / EC."
$d.Content.Find.Execute($oldText1, $true, $false, $false, $false, $false, $true, 1, $false, $oldText1, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Insert the new bullet paragraph right after "Personal data provided in
#    the form on the landing page ..." (still paragraph 9 at this point,
#    since none of the edits above add/remove paragraphs).
# ---------------------------------------------------------------------------
$p9 = $d.Paragraphs(9)
$p9.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(10)
$newPara.Range.Text = "We will process your data in order to: prepare an educational offer, conclude an educational contract, service the educational contract, other activities necessary to perform the contract, send and provide information on the offers offered by the Konwersatorium Muzyczne - by post, e-mail, text message, by phone - where telephone calls will be recorded for which you consent."

# ---------------------------------------------------------------------------
# 4) Merge the runs around " legitimate interests ... " / lastRenderedPageBreak
#    / "for the purpose consistent ..." into a single run, and add a
#    lastRenderedPageBreak before "Personal data is processed: a. ...".
# ---------------------------------------------------------------------------
$oldText2 = " legitimate interests (legitimate purposes), and the processing does not violate the rights and freedoms of the data subject to the extent and for the purpose consistent with the consent expressed by you if you [for example] subscribed to the newsletter."
$d.Content.Find.Execute($oldText2, $true, $false, $false, $false, $false, $true, 1, $false, $oldText2, 2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Merge the split runs around "e.g." into a single run (removes the now
#    redundant proofErr gramStart/gramEnd wrapper). The match text is kept
#    clear of the apostrophe in "User's" so the replace engine's smart-quote
#    autocorrect does not mangle that (untouched) character.
# ---------------------------------------------------------------------------
$oldText3 = "applicable law (e.g. law enforcement authorities)."
$d.Content.Find.Execute($oldText3, $true, $false, $false, $false, $false, $true, 1, $false, $oldText3, 2) | Out-Null
